$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New price / volume(1h) figures for the cryptos list refresh.
# D-column values that look like plain decimals (e.g. "592.91") need a
# leading apostrophe so Excel keeps them as text instead of coercing them
# to a number; the style is then reset to Normal so no numFmt/quotePrefix
# styling sticks to the cell.
$ws.Range("D2").Value = "65.030.02"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "3.517.22"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'592.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "'133.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("D7").Value = "3.515.76"
$ws.Range("E7").Value = "  -1.51%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").Value = "'7.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.08%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "4.117.52"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").Value = "'27.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "3.518.73"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "65.013.46"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("D22").Value = "'392.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "3.659.14"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -4.88%  "
$ws.Range("E28").Value = "  +8.57%  "
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("D33").Value = "3.523.71"
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("D34").Value = "'24.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "'5.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.85%  "
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "'167.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("E43").Value = "  +4.29%  "
$ws.Range("D44").Value = "'42.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").Value = "'25.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.19%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'4.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "'6.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "2.422.21"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("D51").Value = "'0.903"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.05%  "

